$wb = $excel.ActiveWorkbook

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1050.125
$ws.Range("I43").Value = 699.5
$ws.Range("J43").Value = 1167
$ws.Range("K43").Value = 699.5
$ws.Range("L43").Value = 1167
$ws.Range("M43").Value = -630.5
$ws.Range("N43").Value = -1305

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3419.5454
$ws.Range("I76").Value = 3401.6667
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 3401.6667
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -3086.6667
$ws.Range("N76").Value = -4130

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3419.5454
$ws.Range("I79").Value = 3401.6667
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 3401.6667
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -2309.6667
$ws.Range("N79").Value = -5684

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4885.6
$ws.Range("I116").Value = 2698.2856
$ws.Range("K116").Value = 2698.2856
$ws.Range("M116").Value = 743.7143999999998

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1882.8
$ws.Range("I127").Value = 1000
$ws.Range("J127").Value = 2103.5
$ws.Range("K127").Value = 3000
$ws.Range("L127").Value = 6310.5
$ws.Range("M127").Value = 1960
$ws.Range("N127").Value = -16230.5

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2207.5
$ws.Range("I131").Value = 1093.75
$ws.Range("J131").Value = 2950
$ws.Range("K131").Value = 3281.25
$ws.Range("L131").Value = 8850
$ws.Range("M131").Value = 1758.75
$ws.Range("N131").Value = -18930

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7276.6865
$ws.Range("I32").Value = 5140.9614
$ws.Range("K32").Value = 5140.9614
$ws.Range("M32").Value = -4853.9614

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3311.0688
$ws.Range("I61").Value = 3273.3809
$ws.Range("K61").Value = 3273.3809
$ws.Range("M61").Value = -3061.3809

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2908.3333
$ws.Range("J63").Value = 2887.5
$ws.Range("L63").Value = 2887.5
$ws.Range("N63").Value = -4259.5

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2908.3333
$ws.Range("J66").Value = 2887.5
$ws.Range("L66").Value = 14437.5
$ws.Range("N66").Value = -21301.5

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3311.0688
$ws.Range("I136").Value = 3273.3809
$ws.Range("K136").Value = 9820.1427
$ws.Range("M136").Value = -7270.1427

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1864.7727
$ws.Range("I105").Value = 1705.5555
$ws.Range("K105").Value = 1705.5555
$ws.Range("M105").Value = 41.44450000000006

# CRP row 57
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 18000
$ws.Range("J57").Value = 18000
$ws.Range("L57").Value = 18000
$ws.Range("N57").Value = -19120

# CRP row 116
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 34875
$ws.Range("J116").Value = 34875
$ws.Range("L116").Value = 34875
$ws.Range("N116").Value = -44053

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1241.05
$ws.Range("I5").Value = 623.5
$ws.Range("J5").Value = 1652.75
$ws.Range("K5").Value = 1870.5
$ws.Range("L5").Value = 4958.25
$ws.Range("M5").Value = -1758.5
$ws.Range("N5").Value = -5182.25

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 809.75
$ws.Range("I23").Value = 19.5
$ws.Range("J23").Value = 1600
$ws.Range("K23").Value = 58.5
$ws.Range("L23").Value = 4800
$ws.Range("M23").Value = 176.5
$ws.Range("N23").Value = -5270

# CUL row 74
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# CUL row 77
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 735.64
$ws.Range("J131").Value = 735.64
$ws.Range("L131").Value = 2206.92
$ws.Range("N131").Value = -12286.92

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1241.05
$ws.Range("I135").Value = 623.5
$ws.Range("J135").Value = 1652.75
$ws.Range("K135").Value = 5611.5
$ws.Range("L135").Value = 14874.75
$ws.Range("M135").Value = -3076.5
$ws.Range("N135").Value = -19944.75

# CUL row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 3044.6
$ws.Range("I136").Value = 999.5714
$ws.Range("J136").Value = 4834
$ws.Range("K136").Value = 2998.7142
$ws.Range("L136").Value = 14502
$ws.Range("M136").Value = 2101.2858
$ws.Range("N136").Value = -24702

# GSM row 15
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 14997.5
$ws.Range("J15").Value = 14997.5
$ws.Range("L15").Value = 14997.5
$ws.Range("N15").Value = -15573.5

# GSM row 29
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 900
$ws.Range("J29").Value = 900
$ws.Range("L29").Value = 900
$ws.Range("N29").Value = -1480

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16626
$ws.Range("I70").Value = 5666.6665
$ws.Range("J70").Value = 23201.6
$ws.Range("K70").Value = 5666.6665
$ws.Range("L70").Value = 23201.6
$ws.Range("M70").Value = -5396.6665
$ws.Range("N70").Value = -23741.6

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 16626
$ws.Range("I73").Value = 5666.6665
$ws.Range("J73").Value = 23201.6
$ws.Range("K73").Value = 5666.6665
$ws.Range("L73").Value = 23201.6
$ws.Range("M73").Value = -4730.6665
$ws.Range("N73").Value = -25073.6

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3959.375
$ws.Range("I80").Value = 3440
$ws.Range("J80").Value = 4195.4546
$ws.Range("K80").Value = 3440
$ws.Range("L80").Value = 4195.4546
$ws.Range("M80").Value = -2442
$ws.Range("N80").Value = -6191.4546

# GSM row 81
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 14997.5
$ws.Range("J81").Value = 14997.5
$ws.Range("L81").Value = 14997.5
$ws.Range("N81").Value = -16993.5

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3959.375
$ws.Range("I83").Value = 3440
$ws.Range("J83").Value = 4195.4546
$ws.Range("K83").Value = 17200
$ws.Range("L83").Value = 20977.273
$ws.Range("M83").Value = -12208
$ws.Range("N83").Value = -30961.273

# GSM row 84
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 14997.5
$ws.Range("J84").Value = 14997.5
$ws.Range("L84").Value = 44992.5
$ws.Range("N84").Value = -54976.5

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1723.2858
$ws.Range("J136").Value = 1995.8
$ws.Range("L136").Value = 5987.4
$ws.Range("N136").Value = -11087.4
